$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.867.16"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "2.349.23"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.50"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.47"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E8").Value = "  +6.31%  "
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.53"
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("E11").Value = "  -1.60%  "
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "2.768.74"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.70"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").Value = "57.800.52"
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "2.360.32"
$ws.Range("E17").Value = "  +1.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.71"
$ws.Range("E18").Value = "  +1.35%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "329.76"
$ws.Range("E19").Value = "  -2.54%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.28"
$ws.Range("E20").Value = "  +1.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.71"
$ws.Range("E21").Value = "  -2.17%  "
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "62.52"
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("E24").Value = "  -1.85%  "
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.41"
$ws.Range("E26").Value = "  -1.29%  "
$ws.Range("E27").Value = "  -3.05%  "
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.37"
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("D30").Value = "0.0₃0735"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.12"
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.03"
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.36"
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.21"
$ws.Range("E36").Value = "  +1.43%  "
$ws.Range("E37").Value = "  -2.18%  "
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "39.01"
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "142.53"
$ws.Range("E40").Value = "  -4.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.377"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "289.20"
$ws.Range("E43").Value = "  +2.02%  "
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.19"
$ws.Range("E46").Value = "  -0.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.565"
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.385"
$ws.Range("E49").Value = "  +1.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.45"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.07"
$ws.Range("E51").Value = "  +0.44%  "
